$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "68.254.81"
Set-TextValue $ws.Range("E2") "  +0.72%  "
Set-TextValue $ws.Range("D3") "3.349.54"
Set-TextValue $ws.Range("E3") "  +0.48%  "
Set-TextValue $ws.Range("E4") "  +0.06%  "
Set-TextValue $ws.Range("D5") "583.32"
Set-TextValue $ws.Range("E5") "  +0.14%  "
Set-TextValue $ws.Range("D6") "177.38"
Set-TextValue $ws.Range("E6") "  +0.59%  "
Set-TextValue $ws.Range("D7") "1.00"
Set-TextValue $ws.Range("E7") "  +0.13%  "
Set-TextValue $ws.Range("E8") "  +0.25%  "
Set-TextValue $ws.Range("E9") "  +2.96%  "
Set-TextValue $ws.Range("D10") "0.581"
Set-TextValue $ws.Range("E10") "  +0.73%  "
Set-TextValue $ws.Range("D11") "48.16"
Set-TextValue $ws.Range("E12") "  +1.04%  "
Set-TextValue $ws.Range("D13") "688.80"
Set-TextValue $ws.Range("E13") "  +3.77%  "
Set-TextValue $ws.Range("D14") "3.894.25"
Set-TextValue $ws.Range("E14") "  +0.66%  "
Set-TextValue $ws.Range("D15") "8.42"
Set-TextValue $ws.Range("E15") "  +0.02%  "
Set-TextValue $ws.Range("D16") "68.361.21"
Set-TextValue $ws.Range("E16") "  +0.70%  "
Set-TextValue $ws.Range("D17") "0.119"
Set-TextValue $ws.Range("E17") "  +1.17%  "
Set-TextValue $ws.Range("D18") "3.330.79"
Set-TextValue $ws.Range("E18") "  -0.08%  "
Set-TextValue $ws.Range("D19") "17.45"
Set-TextValue $ws.Range("E19") "  -0.16%  "
Set-TextValue $ws.Range("D20") "11.19"
Set-TextValue $ws.Range("E20") "  +2.05%  "
Set-TextValue $ws.Range("D21") "0.895"
Set-TextValue $ws.Range("E21") "  +0.42%  "
Set-TextValue $ws.Range("E22") "  +0.21%  "
Set-TextValue $ws.Range("E23") "  -0.61%  "
Set-TextValue $ws.Range("D24") "99.99"
Set-TextValue $ws.Range("E24") "  +0.54%  "
Set-TextValue $ws.Range("E25") "  +1.35%  "
Set-TextValue $ws.Range("D26") "2.69"
Set-TextValue $ws.Range("E26") "  +0.41%  "
Set-TextValue $ws.Range("D27") "9.54"
Set-TextValue $ws.Range("E27") "  +2.57%  "
Set-TextValue $ws.Range("D28") "32.96"
Set-TextValue $ws.Range("E28") "  -2.04%  "
Set-TextValue $ws.Range("D29") "8.49"
Set-TextValue $ws.Range("E29") "  +0.52%  "
Set-TextValue $ws.Range("D30") "6.93"
Set-TextValue $ws.Range("E30") "  -7.09%  "
Set-TextValue $ws.Range("D31") "561.24"
Set-TextValue $ws.Range("E31") "  -5.26%  "
Set-TextValue $ws.Range("D32") "11.06"
Set-TextValue $ws.Range("E32") "  +0.95%  "
Set-TextValue $ws.Range("D33") "0.105"
Set-TextValue $ws.Range("E33") "  +0.91%  "
Set-TextValue $ws.Range("D34") "57.74"
Set-TextValue $ws.Range("E34") "  +1.51%  "
Set-TextValue $ws.Range("E35") "  +0.09%  "
Set-TextValue $ws.Range("D36") "3.700.01"
Set-TextValue $ws.Range("E36") "  -0.65%  "
Set-TextValue $ws.Range("D37") "3.29"
Set-TextValue $ws.Range("E37") "  +0.31%  "
Set-TextValue $ws.Range("E38") "  +3.96%  "
Set-TextValue $ws.Range("E39") "  +3.18%  "
Set-TextValue $ws.Range("D40") "3.17"
Set-TextValue $ws.Range("E40") "  +1.74%  "
Set-TextValue $ws.Range("D41") "2.61"
Set-TextValue $ws.Range("E41") "  -0.97%  "
Set-TextValue $ws.Range("D42") "0.0₃0672"
Set-TextValue $ws.Range("E42") "  +0.80%  "
Set-TextValue $ws.Range("D43") "0.335"
Set-TextValue $ws.Range("E43") "  +0.41%  "
Set-TextValue $ws.Range("E44") "  +0.74%  "
Set-TextValue $ws.Range("D45") "0.0411"
Set-TextValue $ws.Range("E45") "  +0.94%  "
Set-TextValue $ws.Range("D46") "2.65"
Set-TextValue $ws.Range("E46") "  +2.04%  "
Set-TextValue $ws.Range("E47") "  +0.57%  "
Set-TextValue $ws.Range("E48") "  +0.02%  "
Set-TextValue $ws.Range("E49") "  -0.44%  "
Set-TextValue $ws.Range("D50") "130.61"
Set-TextValue $ws.Range("E50") "  +2.75%  "
Set-TextValue $ws.Range("D51") "2.56"
Set-TextValue $ws.Range("E51") "  -0.76%  "
